$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 25.29000000000051
$ws.Range("H2").Value = (5.921189464667501 * [Math]::Pow(10, -16))
$ws.Range("K2").Value = 38.4061039869934
$ws.Range("L2").Value = "[28.87151546592672, 47.94069250806008]"
$ws.Range("M2").Value = (5.062616992290714 * [Math]::Pow(10, -14))
$ws.Range("N2").Value = (1.012523398458143 * [Math]::Pow(10, -13))
$ws.Range("O2").Value = 1.767342413731195
$ws.Range("P2").Value = "[1.50318447288881, 2.0315003545735806]"
$ws.Range("S2").Value = 58.65628655812257
$ws.Range("T2").Value = "[53.15005027805162, 64.16252283819352]"
$ws.Range("W2").Value = 18.17639639639676
$ws.Range("X2").Value = 17.1131531531535
$ws.Range("Y2").Value = 19.23963963964003

# Row 3
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 22.84000000000013
$ws.Range("G3").Value = (1.474677047141881 * [Math]::Pow(10, -10))
$ws.Range("H3").Value = (43.40052234734469 * [Math]::Pow(10, -11))
$ws.Range("K3").Value = 43.98535132898458
$ws.Range("L3").Value = "[29.204441836134386, 58.76626082183477]"
$ws.Range("M3").Value = (2798.461706809974 * [Math]::Pow(10, -11))
$ws.Range("N3").Value = (2798.461706809974 * [Math]::Pow(10, -11))
$ws.Range("O3").Value = -1.811368737204925
$ws.Range("P3").Value = "[-2.201316173686541, -1.4214213007233099]"
$ws.Range("Q3").Value = (4.440892098500626 * [Math]::Pow(10, -16))
$ws.Range("R3").Value = (4.440892098500626 * [Math]::Pow(10, -16))
$ws.Range("S3").Value = 53.44097468556694
$ws.Range("T3").Value = "[44.93770910568532, 61.94424026544855]"
$ws.Range("W3").Value = 6.58450450450454
$ws.Range("X3").Value = 5.167007007007036
$ws.Range("Y3").Value = 8.002002002002044
